$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns are stored as text in the sheet.
# Values that look like plain decimals (e.g. "1.00", "307.57") would be
# auto-converted to numbers by Excel's smart-entry parsing, so those are
# written with a leading apostrophe (forces text) and then the cell style
# is reset to "Normal" so no stray formatting / quote-prefix marker is left
# behind. Values that already read as text (multi-dot price strings,
# percentage strings with padding spaces, etc.) are set directly.

$ws.Range("D2").Value = '39.859.20'
$ws.Range("E2").Value = '  -4.45%  '
$ws.Range("D3").Value = '2.319.62'
$ws.Range("E3").Value = '  -6.03%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''307.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.20%  '
$ws.Range("D6").Value = '''84.38'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -8.52%  '
$ws.Range("D7").Value = '''0.530'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.54%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = '''0.482'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.92%  '
$ws.Range("D10").Value = '''0.0809'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.43%  '
$ws.Range("D11").Value = '''29.91'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -9.14%  '
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("D13").Value = '2.685.60'
$ws.Range("D14").Value = '''6.38'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.21%  '
$ws.Range("D15").Value = '''14.62'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.63%  '
$ws.Range("D16").Value = '2.323.66'
$ws.Range("E16").Value = '  -5.89%  '
$ws.Range("D17").Value = '''0.751'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.93%  '
$ws.Range("D18").Value = '39.848.11'
$ws.Range("E18").Value = '  -4.32%  '
$ws.Range("D19").Value = '0.0₃0900'
$ws.Range("E19").Value = '  -4.22%  '
$ws.Range("D20").Value = '''6.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.96%  '
$ws.Range("D21").Value = '''67.32'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.33%  '
$ws.Range("D22").Value = '''10.55'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.97%  '
$ws.Range("D23").Value = '''234.56'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.07%  '
$ws.Range("E24").Value = '  -7.33%  '
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("E26").Value = '  -7.03%  '
$ws.Range("D27").Value = '''23.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.64%  '
$ws.Range("E28").Value = '  -4.50%  '
$ws.Range("D29").Value = '''9.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.87%  '
$ws.Range("D30").Value = '''34.88'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.60%  '
$ws.Range("D31").Value = '''151.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.70%  '
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("D33").Value = '''5.06'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.76%  '
$ws.Range("D34").Value = '''2.44'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.88%  '
$ws.Range("D35").Value = '''0.0718'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.91%  '
$ws.Range("E36").Value = '  -2.84%  '
$ws.Range("D37").Value = '''0.0994'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.55%  '
$ws.Range("D38").Value = '''2.73'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.43%  '
$ws.Range("D39").Value = '''15.49'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.03%  '
$ws.Range("E40").Value = '  -7.28%  '
$ws.Range("D41").Value = '''3.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.73%  '
$ws.Range("E42").Value = '  -2.74%  '
$ws.Range("D43").Value = '1.941.49'
$ws.Range("E43").Value = '  -2.99%  '
$ws.Range("E44").Value = '  -6.12%  '
$ws.Range("D45").Value = '''17.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.82%  '
$ws.Range("D46").Value = '''9.32'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.71%  '
$ws.Range("E47").Value = '  -9.74%  '
$ws.Range("D48").Value = '2.551.93'
$ws.Range("E48").Value = '  -6.50%  '
$ws.Range("D49").Value = '''92.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.91%  '
$ws.Range("D50").Value = '''70.26'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.20%  '
$ws.Range("D51").Value = '''50.02'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.14%  '
